$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("b9845_d170212", 9845)
    ,@("b9845_d170215", 9845)
    ,@("b9845_d170216", 9845)
    ,@("b9845_d170217", 9845)
    ,@("b9845_d170219", 9845)
    ,@("b9845_d170220", 9845)
    ,@("b9845_d170221", 9845)
    ,@("b9845_d170222", 9845)
    ,@("b9845_d170223", 9845)
    ,@("b9845_d170516", 9845)
    ,@("b9845_d170517", 9845)
    ,@("b9845_d170518", 9845)
    ,@("b9845_d170525", 9845)
    ,@("b9845_d170527", 9845)
    ,@("b9845_d170528", 9845)
    ,@("b9845_d170529", 9845)
    ,@("b9845_d170603", 9845)
    ,@("b9845_d170605", 9845)
    ,@("b9845_d170606", 9845)
    ,@("b9845_d170612", 9845)
    ,@("b9845_d170614", 9845)
    ,@("b9845_d170622", 9845)
    ,@("b2311_d191218", 2311)
    ,@("b2311_d191219", 2311)
    ,@("b2311_d191220", 2311)
    ,@("b2311_d191222", 2311)
    ,@("b2311_d191223", 2311)
    ,@("b2311_d191224", 2311)
    ,@("b2311_d191225", 2311)
    ,@("b2311_d191226", 2311)
    ,@("b2311_d191229", 2311)
    ,@("b2311_d191230", 2311)
    ,@("b2311_d191231", 2311)
    ,@("b2311_d200101", 2311)
    ,@("b2311_d200102", 2311)
)

$startRow = 210
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $expId = $data[$i][0]
    $batNum = $data[$i][1]
    $ws.Cells.Item($r, 1).Value2 = $expId
    $ws.Cells.Item($r, 2).Value2 = $batNum
    $ws.Cells.Item($r, 3).Value2 = 0
}

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 208
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A245").Select()
